# Add 8 new "Lore of Shadow" spells to the Spells sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Spells")

$rows = @(
    @("Mindrazors",          "New", 5, "Illusion",     "No","Yes","No","No","No","No","No","Yes","No","1.0.0","Complete","Publicly Released","Not on website"),
    @("Mirror Dance",        "New", 2, "Conjuration",  "No","Yes","No","No","No","No","Yes","Yes","Yes","1.0.0","Complete","Publicly Released","Not on website"),
    @("Mystifying Miasma",   "New", 1, "Illusion",     "No","Yes","No","No","No","Yes","No","No","No","1.0.0","Complete","Publicly Released","Not on website"),
    @("Penumbral Pendulum",  "New", 4, "Illusion",     "No","No","No","No","No","No","Yes","No","Yes","2.0.0","Complete","Publicly Released","Not on website"),
    @("Pit of Shadows",      "New", 6, "Abjuration",   "No","No","Yes","No","No","No","Yes","Yes","Yes","1.0.0","Complete","Publicly Released","Not on website"),
    @("Shroud of Dusk",      "New", 1, "Evocation",    "No","No","No","No","No","No","Yes","Yes","Yes","1.1.0","Complete","Publicly Released","Not on website"),
    @("Unseen Lurker",       "New", 3, "Trasmutation", "No","No","No","No","No","Yes","No","No","No","1.0.0","Complete","Publicly Released","Not on website"),
    @("Withering Resolve",   "New", 3, "Enchantment",  "No","Yes","No","No","No","No","No","No","No","1.0.0","Complete","Publicly Released","Not on website")
)

$startRow = 76
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
    $ws.Cells.Item($r, 7).Value = $data[6]
    $ws.Cells.Item($r, 8).Value = $data[7]
    $ws.Cells.Item($r, 9).Value = $data[8]
    $ws.Cells.Item($r, 10).Value = $data[9]
    $ws.Cells.Item($r, 11).Value = $data[10]
    $ws.Cells.Item($r, 12).Value = $data[11]
    $ws.Cells.Item($r, 13).Value = $data[12]
    $ws.Cells.Item($r, 14).Value = $data[13]
    $ws.Cells.Item($r, 15).Value = $data[14]
    $ws.Cells.Item($r, 16).Value = $data[15]
    $ws.Cells.Item($r, 17).Value = $data[16]
}

# Update view: active cell and scroll position
$ws.Application.ActiveWindow.ScrollRow = 59
$ws.Range("A83").Select()
